# Weekly data refresh: insert a new pair of price rows (Primera/Segunda)
# for Acelga @ Terminal La Palmera de La Serena, shifting the existing
# rows 230-269 down to 232-271.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 230:231 (existing rows 230+ shift down to 232+,
# carrying their formatting - including the date style on column D - with them).
$ws.Range("A230:R231").Insert()

# Populate the two newly-inserted rows with this week's data.
$ws.Range("A230").Value = 8
$ws.Range("B230").Value = "Terminal La Palmera de La Serena"
$ws.Range("C230").Value = "Coquimbo"
$ws.Range("D230").Value = 44522
$ws.Range("E230").Value = 4
$ws.Range("F230").Value = 100112009
$ws.Range("G230").Value = "Acelga"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 2400
$ws.Range("K230").Value = 550
$ws.Range("L230").Value = 600
$ws.Range("M230").Value = 575
$ws.Range("N230").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O230").Value = "Provincia del Elquí"
$ws.Range("P230").Value = 288
$ws.Range("Q230").Value = 2
$ws.Range("R230").Value = "Hortaliza"

$ws.Range("A231").Value = 8
$ws.Range("B231").Value = "Terminal La Palmera de La Serena"
$ws.Range("C231").Value = "Coquimbo"
$ws.Range("D231").Value = 44522
$ws.Range("E231").Value = 4
$ws.Range("F231").Value = 100112009
$ws.Range("G231").Value = "Acelga"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Segunda"
$ws.Range("J231").Value = 1400
$ws.Range("K231").Value = 450
$ws.Range("L231").Value = 500
$ws.Range("M231").Value = 475
$ws.Range("N231").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O231").Value = "Provincia del Elquí"
$ws.Range("P231").Value = 238
$ws.Range("Q231").Value = 2
$ws.Range("R231").Value = "Hortaliza"
